$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 213, shifting the existing rows 213:257 down to 215:259.
$ws.Rows("213:214").Insert()

# --- New row 213 ---
$ws.Range("A213").Value = 10
$ws.Range("B213").Value = "Vega Modelo de Temuco"
$ws.Range("C213").Value = "La Araucanía"
$ws.Range("D213").Value = 44511
$ws.Range("E213").Value = 9
$ws.Range("F213").Value = 100112040
$ws.Range("G213").Value = "Cilantro"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 60
$ws.Range("K213").Value = 4000
$ws.Range("L213").Value = 4000
$ws.Range("M213").Value = 4000
$ws.Range("N213").Value = "$/docena de atados (2 kilos)"
$ws.Range("O213").Value = "Provincia de Cautín"
$ws.Range("P213").Value = 2000
$ws.Range("Q213").Value = 2
$ws.Range("R213").Value = "Hortaliza"

# --- New row 214 ---
$ws.Range("A214").Value = 10
$ws.Range("B214").Value = "Vega Modelo de Temuco"
$ws.Range("C214").Value = "La Araucanía"
$ws.Range("D214").Value = 44511
$ws.Range("E214").Value = 9
$ws.Range("F214").Value = 100112040
$ws.Range("G214").Value = "Cilantro"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 30
$ws.Range("K214").Value = 2600
$ws.Range("L214").Value = 2600
$ws.Range("M214").Value = 2600
$ws.Range("N214").Value = "$/docena de atados (2 kilos)"
$ws.Range("O214").Value = "Región del Maule"
$ws.Range("P214").Value = 1300
$ws.Range("Q214").Value = 2
$ws.Range("R214").Value = "Hortaliza"
